$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A33").Value = "Generic"
$ws.Range("B33").Value = "run a cmd at background"
$ws.Range("C33").Value = "Run a cmd at background:`n> {command_body} &`nCheck the process:`n> top`n> ps -a`n> pstree"
$ws.Range("C33").WrapText = $true
